# Updated cryptos list with GitHub Actions run - refresh price/volume figures,
# and re-rank dogwifhat above Bittensor (rows 46/47 swap names/links/values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces plain-text entry for Price values that would otherwise be
# auto-parsed as numbers (e.g. "1.00", "0.511"), matching the inlineStr cells
# in the workbook. Values with two dots (thousands separators, e.g.
# "66.756.61") are never auto-numeric so need no prefix.
$ws.Range("D2").Value = "66.756.61"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.071.12"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'575.56"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'170.29"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.069.56"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").Value = "'0.511"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("E12").Value = "  -3.29%  "
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("D14").Value = "'35.78"
$ws.Range("E14").Value = "  -3.69%  "
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "3.586.10"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "66.765.48"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "'6.99"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "'16.93"
$ws.Range("E19").Value = "  +3.91%  "
$ws.Range("D20").Value = "3.067.42"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").Value = "'490.75"
$ws.Range("E21").Value = "  +3.03%  "
$ws.Range("D22").Value = "'7.70"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "'0.687"
$ws.Range("E23").Value = "  -3.73%  "
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").Value = "'12.66"
$ws.Range("E25").Value = "  -5.50%  "
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("D27").Value = "'10.12"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "'7.79"
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").Value = "'27.54"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "0.0₃0912"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'0.949"
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("D37").Value = "'5.58"
$ws.Range("E37").Value = "  -4.69%  "
$ws.Range("D38").Value = "'46.96"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'0.122"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "'1.96"
$ws.Range("E40").Value = "  -5.12%  "
$ws.Range("D41").Value = "'0.299"
$ws.Range("E41").Value = "  -3.77%  "
$ws.Range("D42").Value = "'8.31"
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("D43").Value = "2.752.79"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").Value = "'0.0345"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "'135.33"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.49"
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'365.88"
$ws.Range("E47").Value = "  -5.30%  "
$ws.Range("D49").Value = "'24.66"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("E51").Value = "  -1.92%  "
